{"js": "const replacements = [\n  [\"454\u00f76=\", \"748\u00f75=\"],\n  [\"396\u00f77=\", \"315\u00f74=\"],\n  [\"777\u00f73=\", \"339\u00f78=\"],\n  [\"923\u00f72=\", \"812\u00f79=\"],\n  [\"682\u00f76=\", \"356\u00f74=\"],\n  [\"935\u00f79=\", \"344\u00f73=\"],\n  [\"277\u00f75=\", \"862\u00f79=\"],\n  [\"336\u00f77=\", \"263\u00f72=\"],\n  [\"297\u00f77=\", \"857\u00f78=\"],\n  [\"726\u00f79=\", \"410\u00f79=\"],\n  [\"555\u00f78=\", \"744\u00f72=\"],\n  [\"984\u00f77=\", \"588\u00f78=\"],\n  [\"833\u00f73=\", \"949\u00f75=\"],\n  [\"504\u00f76=\", \"945\u00f74=\"],\n  [\"935\u00f78=\", \"695\u00f77=\"],\n  [\"488\u00f73=\", \"400\u00f73=\"],\n  [\"160\u00f75=\", \"721\u00f79=\"],\n  [\"579\u00f76=\", \"184\u00f75=\"],\n  [\"161\u00f79=\", \"464\u00f72=\"],\n  [\"417\u00f79=\", \"430\u00f74=\"],\n  [\"866\u00f79=\", \"897\u00f76=\"],\n  [\"650\u00f74=\", \"177\u00f79=\"],\n  [\"718\u00f74=\", \"597\u00f77=\"],\n  [\"851\u00f74=\", \"852\u00f74=\"],\n  [\"807\u00f77=\", \"864\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"454\u00f76=\", \"748\u00f75=\"),\n    @(\"396\u00f77=\", \"315\u00f74=\"),\n    @(\"777\u00f73=\", \"339\u00f78=\"),\n    @(\"923\u00f72=\", \"812\u00f79=\"),\n    @(\"682\u00f76=\", \"356\u00f74=\"),\n    @(\"935\u00f79=\", \"344\u00f73=\"),\n    @(\"277\u00f75=\", \"862\u00f79=\"),\n    @(\"336\u00f77=\", \"263\u00f72=\"),\n    @(\"297\u00f77=\", \"857\u00f78=\"),\n    @(\"726\u00f79=\", \"410\u00f79=\"),\n    @(\"555\u00f78=\", \"744\u00f72=\"),\n    @(\"984\u00f77=\", \"588\u00f78=\"),\n    @(\"833\u00f73=\", \"949\u00f75=\"),\n    @(\"504\u00f76=\", \"945\u00f74=\"),\n    @(\"935\u00f78=\", \"695\u00f77=\"),\n    @(\"488\u00f73=\", \"400\u00f73=\"),\n    @(\"160\u00f75=\", \"721\u00f79=\"),\n    @(\"579\u00f76=\", \"184\u00f75=\"),\n    @(\"161\u00f79=\", \"464\u00f72=\"),\n    @(\"417\u00f79=\", \"430\u00f74=\"),\n    @(\"866\u00f79=\", \"897\u00f76=\"),\n    @(\"650\u00f74=\", \"177\u00f79=\"),\n    @(\"718\u00f74=\", \"597\u00f77=\"),\n    @(\"851\u00f74=\", \"852\u00f74=\"),\n    @(\"807\u00f77=\", \"864\u00f75=\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}"}
